# Edit script for Honor_телефоны.xlsx (sheet "Телефоны")
# - Rewrites the price list rows A2:B108 with the updated catalog
#   (several insertions, reorderings and price changes; row count 91 -> 108)
# - Widens column A from 41 to 43 characters
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (Model name, purchase price) for rows 2..108. Row 1 (header) is unchanged.
$data = @(
    ,("Honor Magic V2 16/512GB (VER-N49) Duos", 37749)
    ,("Honor Magic 6 Pro 12/512GB (BVL-N49) Duos", 25850)
    ,("Honor Magic 7 12/256GB (PTP-N29) Duos", 25850)
    ,("Honor 90 Pro 16/256GB (REP-AN00) Duos", 15700)
    ,("Honor 70 Pro+ 12/256GB (HPB-AN00) Duos", 14450)
    ,("Honor 200 12/512GB (ELI-NX9) Duos", 13950)
    ,("Honor 200 12/256GB (ELI-AN00) Duos", 13250)
    ,("Honor 200 12/256GB (ELI-NX9) Duos", 13250)
    ,("Honor 70 Pro+ 8/256GB (HPB-AN00) Duos", 13250)
    ,("Honor 200 8/256GB (ELI-NX9) Duos", 12700)
    ,("Honor 90 12/256GB (REA-NX9) Duos", 11950)
    ,("Honor 90 12/512GB (REA-NX9) Duos", 11950)
    ,("Honor X9c 12/256GB (BRP-NX1) Duos", 11950)
    ,("Honor 70 12/256GB (FNE-AN00) Duos", 10900)
    ,("Honor X9c 12/256GB (BRC-NX1) Duos", 10900)
    ,("Honor 90 8/256GB (REA-NX9) Duos", 9600)
    ,("Honor X60 Pro 8/256GB (BRP-AN00) Duos", 9050)
    ,("Honor X9c 8/256GB (BRC-NX1) Duos", 9050)
    ,("Honor X9c 8/256GB (BRP-NX1) Duos", 9050)
    ,("Honor X9b 12/256GB (ALI-NX1) Duos", 8550)
    ,("Honor 50 8/256GB (NTH-AN00) Duos", 8200)
    ,("Honor 50 8/256GB (NTH-NX9) Duos", 8200)
    ,("Honor X7c 8/512GB (ALT-LX1) Duos", 8200)
    ,("Honor X8c 8/512GB (ABR-LX2) Duos", 7200)
    ,("Honor X9a 8/256GB (RMO-NX1) Duos", 6750)
    ,("Honor X9b 8/256GB (ALI-NX1) Duos", 6750)
    ,("Honor 200 Lite 8/256GB (LLY-NX1) Duos", 6400)
    ,("Honor 50 6/128GB (NTH-NX9) Duos", 6400)
    ,("Honor 50 8/128GB (NTH-NX9) Duos", 6400)
    ,("Honor 70 8/128GB (FNE-NX9) Duos", 6400)
    ,("Honor X8c 8/256GB (ABR-LX1) Duos", 6400)
    ,("Honor X8b 8/256GB (LLY-LX1) Duos", 6050)
    ,("Honor X9a 6/128GB (RMO-NX1) Duos", 6050)
    ,("Honor X8b 8/128GB (LLY-LX1) Duos", 5700)
    ,("Honor X9 6/128GB (ANY-LX1) Duos", 5700)
    ,("Honor X9 6/128GB (ANY-LX2) Duos", 5700)
    ,("Honor 30 8/128GB (BMH-AN10) Duos", 5300)
    ,("Honor 400 Lite 8/256GB (ABR-NX1) Duos", 5300)
    ,("Honor 90 Lite 8/256GB (CRT-NX1) Duos", 5300)
    ,("Honor X8c 6/128GB (ABR-LX1) Duos", 5300)
    ,("Honor X8c 8/128GB (ABR-LX1) Duos", 5300)
    ,("Honor 20 Pro 8/256GB (YAL-L41) Duos", 4950)
    ,("Honor 30S 6/128GB (CDY-NX9A) Duos", 4700)
    ,("Honor X6c 6/256GB (NIC-LX1) Duos", 4300)
    ,("Honor 50 Lite 6/128GB (NTN-LX1) Duos", 4000)
    ,("Honor X7a Plus 6/128GB (RKY-LX1) Duos", 4000)
    ,("Honor X7c 6/128GB (ALT-LX1) Duos", 4000)
    ,("Honor X8a 6/128GB (CRT-LX1) Duos", 4000)
    ,("Honor 20 6/128GB (YAL-L21) Duos", 3700)
    ,("Honor 30i 4/128GB (LRA-LX1) Duos", 3700)
    ,("Honor View 20 6/128GB (PCT-L29) Duos", 3700)
    ,("Honor X6c 6/128GB (NIC-LX1) Duos", 3700)
    ,("Honor X7b 8/128GB (CLK-LX1) Duos", 3700)
    ,("Honor 10 Lite 6/128GB (HRY-LX1) Duos", 3400)
    ,("Honor 10i 6/128GB (HRY-LX1T) Duos", 3400)
    ,("Honor 9X Premium 6/128GB (STK-LX1) Duos", 3400)
    ,("Honor View 10 6/128GB (BKL-L09) Duos", 3400)
    ,("Honor X7a 4/128GB (RKY-LX1) Duos", 3400)
    ,("Honor 10 4/128GB (COL-L29) Duos", 3100)
    ,("Honor 10 6/64GB (COL-AL10) Duos", 3100)
    ,("Honor 20S/Lite 6/128GB (MAR-LX1H) Duos", 3100)
    ,("Honor 9X 4/128GB (STK-LX1) Duos", 3100)
    ,("Honor X7 4/128GB (CMA-LX1) Duos", 3100)
    ,("Honor 10X Lite 4/128GB (DNN-LX9) Duos", 2850)
    ,("Honor 20 Lite 4/128GB (MAR-LX1H) Duos", 2850)
    ,("Honor 8 Pro 6/64GB (DUK-L09) Duos", 2850)
    ,("Honor X6a 4/128GB (WDY-LX1) Duos", 2850)
    ,("Honor 10i 4/128GB (HRY-LX1T) Duos", 2550)
    ,("Honor 8X 4/128GB (JSN-L21) Duos", 2550)
    ,("Honor 10 4/64GB (COL-L29) Duos", 2400)
    ,("Honor 9C 4/64GB (AKA-L29) Duos", 2400)
    ,("Honor X6 4/64GB (VNE-LX1) Duos", 2250)
    ,("Honor 8X 4/64GB (JSN-L21) Duos", 2100)
    ,("Honor 7X 4/64GB (BND-L21) Duos", 1950)
    ,("Honor 10 Lite 3/32GB (HRY-LX1) Duos", 1800)
    ,("Honor 10 Lite 3/64GB (HRY-LX1) Duos", 1800)
    ,("Honor 8C 3/32GB (BKK-L21) Duos", 1800)
    ,("Honor 8S Prime 3/64GB (KSA-LX9) Duos", 1800)
    ,("Honor 9 4/64GB (STF-L09) Duos", 1800)
    ,("Honor 9 Lite 4/64GB (LLD-L31/AL10) Duos", 1800)
    ,("Honor 9A 3/64GB (MOA-LX9N) Duos", 1800)
    ,("Honor 8 4/32GB (FRD-L09) Duos", 1550)
    ,("Honor 8A Prime 3/64GB (JAT-LX1) Duos", 1550)
    ,("Honor 9 Lite 3/32GB (LLD-L31/AL10) Duos", 1550)
    ,("Honor 7 3/16GB (PLK-L01) Duos", 1400)
    ,("Honor X5 2/32GB (VNA-LX2) Duos", 1400)
    ,("Honor 7A Prime 2/32GB (DUA-L22) Duos", 1250)
    ,("Honor 7C 3/32GB (AUM-L41) Duos", 1250)
    ,("Honor 8 Lite 3/32GB (PRA-TL10) Duos", 1250)
    ,("Honor 8 Lite 4/32GB (PRA-TL10) Duos", 1250)
    ,("Honor 8A 2/32GB (JAT-LX1) Duos", 1250)
    ,("Honor 8S 2/32GB (KSA-LX9) Duos", 1250)
    ,("Honor 6C Pro 3/32GB (JMM-L22) Duos", 1100)
    ,("Honor 9S 2/32GB (DUA-LX9) Duos", 1100)
    ,("Honor 5X 2/16GB (KIW-L21) Duos", 900)
    ,("Honor 7A Pro 2/16GB (AUM-L29) Duos", 900)
    ,("Honor 5C 2/16GB (NEM-L51) Duos", 850)
    ,("Honor 6C 3/32GB (DIG-L21HN) Duos", 850)
    ,("Honor 4X 2/8GB (CHE2-L11) Duos", 700)
    ,("Honor 7A 2/16GB (DUA-L22) Duos", 700)
    ,("Honor 7S 1/16GB (DRA-LX5) Duos", 700)
    ,("Honor 7S 2/16GB (DRA-LX5) Duos", 700)
    ,("Honor 4C 2/8GB (CHM-U01) Duos", 550)
    ,("Honor 5A 2/16GB (LYO-L21) Duos", 550)
    ,("Honor 6A 2/16GB (DLI-TL20) Duos", 550)
    ,("Honor 3 (HN3-U00)", 450)
    ,("Honor 3C Lite (HOL-U19) Duos", 450)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# Widen column A from 41 to 43 characters. Excels ColumnWidth setter bakes in
# the standard +5/6 character padding offset, so request 42.15 to land on 43.
$ws.Range("A1").ColumnWidth = 42.15
